$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Long / unique whole-sentence & phrase replacements first ---
# (doing these before the short generic fragments like "and "/" or " avoids
#  those fragments being mis-applied inside sentences that are translated as a whole)

# Language list header / link text
Replace-All "English" "Inglés"
Replace-All " / Portuguese / French / Thai / Vietnamese / Spanish" " / Portugués / Francés / Tailandés / Vietnamita / Español"

# Brief table
Replace-All "Brief" "Breve"
Replace-All "An email sent to partners in the target country who have RSVPed yes. We want them to submit their documents. It will be sent via customer.io" "Un correo electrónico enviado a los socios del país objetivo que han confirmado su asistencia. Queremos que envíen sus documentos. Se enviará a través de customer.io"
Replace-All "Target audience" "Público objetivo"
Replace-All "Invited partners who RSVP yes" "Socios invitados que han confirmado su asistencia"

# Subject line
Replace-All "Subject: " "Asunto: "
Replace-All " — take the next step" " — da el siguiente paso"

# Heading / greeting
Replace-All "Thank you for registering for " "Gracias por inscribirte al "
Replace-All "Hi " "Hola "
Replace-All "We are excited for you to join us at " "¡Estamos encantados de que te unas a nosotros en "
Replace-All "To confirm your registration, we would require you and one guest of your choice to provide us with:" "Para confirmar tu inscripción, necesitamos que tú y un acompañante de tu elección nos proporcionen:"

# Document list
Replace-All "A signed copy of the " "Una copia firmada del "
Replace-All "Code of Conduct " "Código de Conducta "
Replace-All "Terms and Conditions" "Términos y Condiciones"
Replace-All " (1 set from each person)" " (1 por persona)"
Replace-All "A scanned copy of your international passports" "Una copia escaneada de sus pasaportes internacionales"
Replace-All "Covid-19 vaccination certificates" "Certificados de vacunación Covid-19"

# Button
Replace-All "Send my details" "Enviar mis detalles"

# Country manager / event package (long sentences -- consume "and "/" or " occurrences
# that are NOT meant to become standalone "y "/" o ")
Replace-All "Your country manager will be in touch to confirm your booking or request any other relevant details. " "El gestor de tu país se pondrá en contacto contigo para confirmar tu reserva o solicitarte cualquier otro detalle relevante. "
Replace-All "Our event package offers you and your guest: " "Nuestro paquete para eventos te ofrece a ti y a tus invitados: "
Replace-All "Flight tickets " "Pasajes de avión "
Replace-All "Travel insurance " "Seguro de viaje "
Replace-All "Airport – Hotel – Airport transfer " "Aeropuerto - Hotel - Traslado del aeropuerto "
Replace-All "One hotel room for you and your guest / Two hotel rooms for you and your guest" "Una habitación de hotel para ti y tu invitado / Dos habitaciones de hotel para ti y tu invitado"

# Check-in / Check-out
Replace-All "Check-in" "Ingreso"
Replace-All "Check-out" "Salida"
Replace-All " on " " en "
Replace-All "[DD Mmm YYYY]" "[DD Mmm AAAA]"

# Meals / sightseeing
Replace-All "Meals (Breakfast, lunch, and dinner)" "Comidas (desayuno, almuerzo y cena)"
Replace-All "Sightseeing tour of " "Visita turística de "
Replace-All "We will send you a confirmation letter before your departure date with the event agenda and information about your flights, transportation, and accommodation. " "Te enviaremos una carta de confirmación antes de la fecha de salida con el programa del evento e información sobre tus vuelos, transporte y alojamiento. "

# Contact info
Replace-All "If you have any questions, please contact us via " "Si tienes alguna pregunta, entra en contacto con nosotros por "
Replace-All "If you have any questions, please contact your country manager, " "Si tienes alguna pregunta, entra en contacto con el gestor de tu país "
Replace-All ", at " ", en "

# Closing
Replace-All "We look forward to seeing you soon." "Esperamos verte pronto."

# --- Short generic fragments last (now unambiguous, since the surrounding
#     sentences that were NOT supposed to get this exact fragment translated
#     have already been fully replaced above) ---
Replace-All "and " "y "
Replace-All " or " " o "

# --- Comments (separate part, not reached by $d.Content.Find) ---
foreach ($c in $d.Comments) {
    switch ($c.Range.Text) {
        "link to T&C" { $c.Range.Text = "enlace a T&C" }
        "link to COC" { $c.Range.Text = "enlace a COC" }
        "please confirm these" { $c.Range.Text = "por favor confirme esto" }
        "choose either one" { $c.Range.Text = "elija uno de los dos" }
        "please check if these are all the required documents" { $c.Range.Text = "por favor verifique si estos son todos los documentos requeridos" }
    }
}
